# Generate Report for Handoff
# - Removes the "e6407dfc-9ddc-40c3-98bc-a15128ca1557.*" row (row 3) from every
#   sheet, so the ".localization-config" row shifts up from row 4 to row 3.
# - Marks the "c4e5f570-3515-4692-bba9-843b4e4ef178.*" entry as freshly handed
#   off again ("Ready for handoff") instead of "Handed back: in sync with en-US".
# - Bumps the "Latest Handoff Datetime" for that entry on both locale sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Range("B2").Value = "Ready for handoff"
$ws1.Range("C2").Value = "Ready for handoff"

# Drop the e6407dfc row entirely; .localization-config shifts from row 4 to row 3.
$ws1.Rows.Item(3).Delete()

# Rebuild the hyperlinks collection to drop stale refs and fix the moved row.
$ws1.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Cells.Item(2,1), "https://github.com/OpenLocalizationTest/oltest/blob/78904320ef8b08d0285accb9596a7141aa86631d/e2e/c4e5f570-3515-4692-bba9-843b4e4ef178.md", [Type]::Missing, [Type]::Missing, "c4e5f570-3515-4692-bba9-843b4e4ef178.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Cells.Item(3,1), "https://github.com/OpenLocalizationTest/oltest/blob/78904320ef8b08d0285accb9596a7141aa86631d/.localization-config", [Type]::Missing, [Type]::Missing, ".localization-config") | Out-Null

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Range("B2").Value = "Ready for handoff"
$ws2.Range("D2").Value = "2016-03-09 07:07:06"

$ws2.Rows.Item(3).Delete()

$ws2.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Cells.Item(2,1), "https://github.com/OpenLocalizationTest/oltest/blob/78904320ef8b08d0285accb9596a7141aa86631d/e2e/c4e5f570-3515-4692-bba9-843b4e4ef178.md", [Type]::Missing, [Type]::Missing, "c4e5f570-3515-4692-bba9-843b4e4ef178.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Cells.Item(2,3), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8487069c711053a555c3858caf5dac8b07aa77b7/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/c4e5f570-3515-4692-bba9-843b4e4ef178.e5662deb35a34c1dde18586d8401bcc4645dbfc9.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "c4e5f570-3515-4692-bba9-843b4e4ef178.e5662deb35a34c1dde18586d8401bcc4645dbfc9.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Cells.Item(2,5), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/35a6e12972cb7831802f7ca498ba43f18234a0c7/e2e/c4e5f570-3515-4692-bba9-843b4e4ef178.md", [Type]::Missing, [Type]::Missing, "c4e5f570-3515-4692-bba9-843b4e4ef178.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Cells.Item(2,6), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/6286ea4534d8036bc38b26739c7e7e08e18c36d9/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/c4e5f570-3515-4692-bba9-843b4e4ef178.e5662deb35a34c1dde18586d8401bcc4645dbfc9.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "c4e5f570-3515-4692-bba9-843b4e4ef178.e5662deb35a34c1dde18586d8401bcc4645dbfc9.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Cells.Item(3,1), "https://github.com/OpenLocalizationTest/oltest/blob/78904320ef8b08d0285accb9596a7141aa86631d/.localization-config", [Type]::Missing, [Type]::Missing, ".localization-config") | Out-Null

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Range("B2").Value = "Ready for handoff"
$ws3.Range("D2").Value = "2016-03-09 07:07:10"

$ws3.Rows.Item(3).Delete()

$ws3.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Cells.Item(2,1), "https://github.com/OpenLocalizationTest/oltest/blob/78904320ef8b08d0285accb9596a7141aa86631d/e2e/c4e5f570-3515-4692-bba9-843b4e4ef178.md", [Type]::Missing, [Type]::Missing, "c4e5f570-3515-4692-bba9-843b4e4ef178.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Cells.Item(2,3), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d8fe1182bfab1aeb66c0792b9b47ce61dcd387cf/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/c4e5f570-3515-4692-bba9-843b4e4ef178.e5662deb35a34c1dde18586d8401bcc4645dbfc9.de-de.xlf", [Type]::Missing, [Type]::Missing, "c4e5f570-3515-4692-bba9-843b4e4ef178.e5662deb35a34c1dde18586d8401bcc4645dbfc9.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Cells.Item(2,5), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/c61e2f40474c7bd143a781ce040cb8a3a1f3378a/e2e/c4e5f570-3515-4692-bba9-843b4e4ef178.md", [Type]::Missing, [Type]::Missing, "c4e5f570-3515-4692-bba9-843b4e4ef178.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Cells.Item(2,6), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/75a71751589415e0f36aab6cf8e7f1bfa067e925/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/c4e5f570-3515-4692-bba9-843b4e4ef178.e5662deb35a34c1dde18586d8401bcc4645dbfc9.de-de.xlf", [Type]::Missing, [Type]::Missing, "c4e5f570-3515-4692-bba9-843b4e4ef178.e5662deb35a34c1dde18586d8401bcc4645dbfc9.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Cells.Item(3,1), "https://github.com/OpenLocalizationTest/oltest/blob/78904320ef8b08d0285accb9596a7141aa86631d/.localization-config", [Type]::Missing, [Type]::Missing, ".localization-config") | Out-Null
